$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.027209544566925
$ws.Range("D2").Value = 1.036216243793162
$ws.Range("E2").Value = 1.030840397174019
$ws.Range("F2").Value = 1.045018289888686
$ws.Range("I2").Value = 1.035126635783653
$ws.Range("J2").Value = 1.03236868986662
$ws.Range("K2").Value = 1.039010620174424
$ws.Range("L2").Value = 1.0336502439994
$ws.Range("M2").Value = 1.047787708967846
$ws.Range("N2").Value = 1.01480623656739
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.028123558461283
$ws.Range("D3").Value = 1.036907023967284
$ws.Range("E3").Value = 1.031701872351562
$ws.Range("F3").Value = 1.045848965104229
$ws.Range("I3").Value = 1.035301004861554
$ws.Range("J3").Value = 1.032922973981904
$ws.Range("K3").Value = 1.039511271360815
$ws.Range("L3").Value = 1.034320015724604
$ws.Range("M3").Value = 1.048429689325349
$ws.Range("N3").Value = 1.01499076227881
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.028715397392193
$ws.Range("D4").Value = 1.037354013199708
$ws.Range("E4").Value = 1.032260065816473
$ws.Range("F4").Value = 1.046386795153065
$ws.Range("I4").Value = 1.035412178216846
$ws.Range("J4").Value = 1.033281429701414
$ws.Range("K4").Value = 1.039834558552646
$ws.Range("L4").Value = 1.034753525362923
$ws.Range("M4").Value = 1.04884477341377
$ws.Range("N4").Value = 1.01511005615366
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.0289643036803
$ws.Range("D5").Value = 1.037541927871446
$ws.Range("E5").Value = 1.032494911152504
$ws.Range("F5").Value = 1.046612975677337
$ws.Range("I5").Value = 1.035458518486731
$ws.Range("J5").Value = 1.03343207477131
$ws.Range("K5").Value = 1.039970307457184
$ws.Range("L5").Value = 1.034935800915261
$ws.Range("M5").Value = 1.049019196711796
$ws.Range("N5").Value = 1.015160181379892
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.029006101807852
$ws.Range("D6").Value = 1.037573479548989
$ws.Range("E6").Value = 1.032534353277067
$ws.Range("F6").Value = 1.046650956836029
$ws.Range("I6").Value = 1.035466275918587
$ws.Range("J6").Value = 1.033457365794807
$ws.Range("K6").Value = 1.039993090817534
$ws.Range("L6").Value = 1.03496640739436
$ws.Range("M6").Value = 1.049048478532831
$ws.Range("N6").Value = 1.015168596094442
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.0287187229097
$ws.Range("D7").Value = 1.037356524125181
$ws.Range("E7").Value = 1.032263203122778
$ws.Range("F7").Value = 1.046389817088741
$ws.Range("I7").Value = 1.03541279897875
$ws.Range("J7").Value = 1.03328344282491
$ws.Range("K7").Value = 1.039836373070416
$ws.Range("L7").Value = 1.034755960828917
$ws.Range("M7").Value = 1.048847104374893
$ws.Range("N7").Value = 1.015110726031732
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.027518354086529
$ws.Range("D8").Value = 1.036449693526179
$ws.Range("E8").Value = 1.031131378196275
$ws.Range("F8").Value = 1.045298951627158
$ws.Range("I8").Value = 1.035185906863898
$ws.Range("J8").Value = 1.032556054199371
$ws.Range("K8").Value = 1.039179955095141
$ws.Range("L8").Value = 1.033876570372661
$ws.Range("M8").Value = 1.048004734672517
$ws.Range("N8").Value = 1.014868619728209
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.025406347948763
$ws.Range("D9").Value = 1.034851875057427
$ws.Range("E9").Value = 1.029142854951191
$ws.Range("F9").Value = 1.04337928996272
$ws.Range("I9").Value = 1.034773451936769
$ws.Range("J9").Value = 1.031272793854434
$ws.Range("K9").Value = 1.038018202720693
$ws.Range("L9").Value = 1.032327965082111
$ws.Range("M9").Value = 1.04651798197429
$ws.Range("N9").Value = 1.014441198616276
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.02400056047471
$ws.Range("D10").Value = 1.033786845206175
$ws.Range("E10").Value = 1.027821224779178
$ws.Range("F10").Value = 1.042101349974221
$ws.Range("I10").Value = 1.034490027556627
$ws.Range("J10").Value = 1.03041633576696
$ws.Range("K10").Value = 1.037240369802547
$ws.Range("L10").Value = 1.031296301941927
$ws.Range("M10").Value = 1.045525289850987
$ws.Range("N10").Value = 1.014155737530591
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023392378468041
$ws.Range("D11").Value = 1.033325738217258
$ws.Range("E11").Value = 1.027249923556701
$ws.Range("F11").Value = 1.041548443371235
$ws.Range("I11").Value = 1.034365305774794
$ws.Range("J11").Value = 1.03004526743587
$ws.Range("K11").Value = 1.036902784101023
$ws.Range("L11").Value = 1.030849771118941
$ws.Range("M11").Value = 1.045095099065471
$ws.Range("N11").Value = 1.014032012762718
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023166553763919
$ws.Range("D12").Value = 1.03315447280044
$ws.Range("E12").Value = 1.027037864522406
$ws.Range("F12").Value = 1.04134313844307
$ws.Range("I12").Value = 1.03431867921962
$ws.Range("J12").Value = 1.029907404659691
$ws.Range("K12").Value = 1.036777273867493
$ws.Range("L12").Value = 1.030683938841192
$ws.Range("M12").Value = 1.044935256046611
$ws.Range("N12").Value = 1.013986038494653
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.023214990239301
$ws.Range("D13").Value = 1.033191209325215
$ws.Range("E13").Value = 1.027083345199322
$ws.Range("F13").Value = 1.041387173889596
$ws.Range("I13").Value = 1.034328694307894
$ws.Range("J13").Value = 1.029936978101101
$ws.Range("K13").Value = 1.036804201463706
$ws.Range("L13").Value = 1.030719509089461
$ws.Range("M13").Value = 1.04496954520316
$ws.Range("N13").Value = 1.013995900911578
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023373710065531
$ws.Range("D14").Value = 1.033311581146494
$ws.Range("E14").Value = 1.027232391669888
$ws.Range("F14").Value = 1.041531471380778
$ws.Range("I14").Value = 1.034361457711341
$ws.Range("J14").Value = 1.03003387229197
$ws.Range("K14").Value = 1.036892411743101
$ws.Range("L14").Value = 1.030836062769578
$ws.Range("M14").Value = 1.045081887424353
$ws.Range("N14").Value = 1.014028212867722
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023471513376713
$ws.Range("D15").Value = 1.033385747602543
$ws.Range("E15").Value = 1.027324243731227
$ws.Range("F15").Value = 1.041620387045829
$ws.Range("I15").Value = 1.034381604682548
$ws.Range("J15").Value = 1.030093567860593
$ws.Range("K15").Value = 1.036946745687594
$ws.Range("L15").Value = 1.030907879237059
$ws.Range("M15").Value = 1.045151098456246
$ws.Range("N15").Value = 1.014048119038694
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02404093483667
$ws.Range("D16").Value = 1.033817448711328
$ws.Range("E16").Value = 1.027859160790396
$ws.Range("F16").Value = 1.042138054180452
$ws.Range("I16").Value = 1.034498262914091
$ws.Range("J16").Value = 1.030440957879737
$ws.Range("K16").Value = 1.037262757941161
$ws.Range("L16").Value = 1.031325940744966
$ws.Range("M16").Value = 1.045553832969188
$ws.Range("N16").Value = 1.014163946277605
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024398261386348
$ws.Range("D17").Value = 1.034088260019899
$ws.Range("E17").Value = 1.028194961743581
$ws.Range("F17").Value = 1.042462894677458
$ws.Range("I17").Value = 1.034570905353076
$ws.Range("J17").Value = 1.030658809267706
$ws.Range("K17").Value = 1.037460776354337
$ws.Range("L17").Value = 1.031588230354374
$ws.Range("M17").Value = 1.045806365333561
$ws.Range("N17").Value = 1.014236570260436
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024606735296788
$ws.Range("D18").Value = 1.034246225164134
$ws.Range("E18").Value = 1.028390922577274
$ws.Range("F18").Value = 1.042652411938753
$ws.Range("I18").Value = 1.034613083714234
$ws.Range("J18").Value = 1.030785857355491
$ws.Range("K18").Value = 1.037576201913594
$ws.Range("L18").Value = 1.031741237324685
$ws.Range("M18").Value = 1.045953629468757
$ws.Range("N18").Value = 1.014278919167425
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024677828195482
$ws.Range("D19").Value = 1.034300088080319
$ws.Range("E19").Value = 1.028457756026969
$ws.Range("F19").Value = 1.042717039721621
$ws.Range("I19").Value = 1.034627432730721
$ws.Range("J19").Value = 1.030829173887043
$ws.Range("K19").Value = 1.037615546250918
$ws.Range("L19").Value = 1.031793411743862
$ws.Range("M19").Value = 1.046003836972266
$ws.Range("N19").Value = 1.014293357099333
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024359918307165
$ws.Range("D20").Value = 1.034059203950438
$ws.Range("E20").Value = 1.028158923764726
$ws.Range("F20").Value = 1.04242803787644
$ws.Range("I20").Value = 1.034563131438115
$ws.Range("J20").Value = 1.030635438034109
$ws.Range("K20").Value = 1.037439538625784
$ws.Range("L20").Value = 1.031560087303676
$ws.Range("M20").Value = 1.04577927447883
$ws.Range("N20").Value = 1.014228779574232
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02332696878343
$ws.Range("D21").Value = 1.033276134350025
$ws.Range("E21").Value = 1.02718849711103
$ws.Range("F21").Value = 1.041488977442847
$ws.Range("I21").Value = 1.034351817962913
$ws.Range("J21").Value = 1.030005340228461
$ws.Range("K21").Value = 1.036866439203844
$ws.Range("L21").Value = 1.030801739801542
$ws.Range("M21").Value = 1.045048806837419
$ws.Range("N21").Value = 1.014018698278706
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022677982584273
$ws.Range("D22").Value = 1.032783847472811
$ws.Range("E22").Value = 1.026579206482379
$ws.Range("F22").Value = 1.040898954058594
$ws.Range("I22").Value = 1.034217225304126
$ws.Range("J22").Value = 1.029608991123944
$ws.Range("K22").Value = 1.036505439012596
$ws.Range("L22").Value = 1.030325105878035
$ws.Range("M22").Value = 1.044589238398681
$ws.Range("N22").Value = 1.013886511461236
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023021977527718
$ws.Range("D23").Value = 1.033044811854517
$ws.Range("E23").Value = 1.02690212143649
$ws.Range("F23").Value = 1.04121169805157
$ws.Range("I23").Value = 1.034288739305298
$ws.Range("J23").Value = 1.029819120155879
$ws.Range("K23").Value = 1.036696875154113
$ws.Range("L23").Value = 1.030577762099766
$ws.Range("M23").Value = 1.044832891780415
$ws.Range("N23").Value = 1.013956595580411
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024377243735509
$ws.Range("D24").Value = 1.034072333120724
$ws.Range("E24").Value = 1.028175207486395
$ws.Range("F24").Value = 1.042443788030817
$ws.Range("I24").Value = 1.034566644731272
$ws.Range("J24").Value = 1.030645998553368
$ws.Range("K24").Value = 1.03744913527371
$ws.Range("L24").Value = 1.031572803881755
$ws.Range("M24").Value = 1.045791515774293
$ws.Range("N24").Value = 1.014232299885347
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.02595196638047
$ws.Range("D25").Value = 1.035264924066986
$ws.Range("E25").Value = 1.029656228063384
$ws.Range("F25").Value = 1.043875252341676
$ws.Range("I25").Value = 1.034881574883464
$ws.Range("J25").Value = 1.031604719681855
$ws.Range("K25").Value = 1.038319135892192
$ws.Range("L25").Value = 1.032728191487067
$ws.Range("M25").Value = 1.046902617801189
$ws.Range("N25").Value = 1.014551789470691
